# Generate Report for Handback
# Update the two source-file UUID placeholders and associated xliff / timestamp
# metadata that a fresh handback run produced.
#
#   old "171df3e6-c226-4652-a681-bb125c39f058" -> new "679d2c86-1021-44b5-97c3-1e3b8ea53ffb"
#   old "d9077830-64b5-469f-b80c-d17bb6746bb1"  -> new "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview.Range("A2").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md"
$overview.Range("B2").Value = "e2e\679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md"
$overview.Range("G2").Value = "2016-08-21 01:07:02"

$overview.Range("A3").Value = "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md"
$overview.Range("B3").Value = "e2e\ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md"
$overview.Range("G3").Value = "2016-08-21 01:07:02"

# refresh the displayed hyperlink text to match the new file names
$overview.Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md", "", "", "e2e\679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md") | Out-Null
$overview.Hyperlinks.Add($overview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/d9077830-64b5-469f-b80c-d17bb6746bb1.md", "", "", "e2e\ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn.Range("A2").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md"
$zhcn.Range("G2").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-08-21 01:06:56"
$zhcn.Range("I2").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md"
$zhcn.Range("J2").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-21 01:07:25"

$zhcn.Range("A3").Value = "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md"
$zhcn.Range("G3").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-21 01:06:56"
$zhcn.Range("I3").Value = "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md"
$zhcn.Range("J3").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-21 01:07:25"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md", "", "", "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c302ffa31e80f4cff226a14d3c49195caa1153ee/e2e/679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md", "", "", "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/d9077830-64b5-469f-b80c-d17bb6746bb1.md", "", "", "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c302ffa31e80f4cff226a14d3c49195caa1153ee/e2e/d9077830-64b5-469f-b80c-d17bb6746bb1.md", "", "", "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md") | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede.Range("A2").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md"
$dede.Range("G2").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.de-de.xlf"
$dede.Range("H2").Value = "2016-08-21 01:07:02"
$dede.Range("I2").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md"
$dede.Range("J2").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.de-de.xlf"
$dede.Range("K2").Value = "2016-08-21 01:07:31"

$dede.Range("A3").Value = "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md"
$dede.Range("G3").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.de-de.xlf"
$dede.Range("H3").Value = "2016-08-21 01:07:02"
$dede.Range("I3").Value = "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md"
$dede.Range("J3").Value = "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.e67babab4b5d0461643d5020868d35f6ef856882.de-de.xlf"
$dede.Range("K3").Value = "2016-08-21 01:07:31"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md", "", "", "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/972acf5a40df5d95fc0d34c680a7629d1a93a5f7/e2e/679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md", "", "", "679d2c86-1021-44b5-97c3-1e3b8ea53ffb.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b47edfe91dd2779a88e5ce69427a492e9740e01/e2e/d9077830-64b5-469f-b80c-d17bb6746bb1.md", "", "", "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/972acf5a40df5d95fc0d34c680a7629d1a93a5f7/e2e/d9077830-64b5-469f-b80c-d17bb6746bb1.md", "", "", "ffff340c7b68-2ecc-4134-83f6-747e86c2b4a2.md") | Out-Null
